# Fix eliminacion de registros previos Excel
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vsSAE_Organizacion")

# Update column A (rows 2-18) from 1 to 8
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = 8
}

# Update the view/selection to A2:A18 (also resets the scrolled topLeftCell)
$ws.Activate()
$ws.Range("A2:A18").Select()
